$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 734.3
$ws.Range("I2").Value = 605.375
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 605.375
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -492.375
$ws.Range("N2").Value = -1476

$ws.Range("H9").Value = 6170.3887
$ws.Range("I9").Value = 10299.4
$ws.Range("J9").Value = 1009.125
$ws.Range("K9").Value = 10299.4
$ws.Range("L9").Value = 1009.125
$ws.Range("M9").Value = -10130.4
$ws.Range("N9").Value = -1347.125

$ws.Range("H33").Value = 299.8
$ws.Range("I33").Value = 296.44446
$ws.Range("J33").Value = 330
$ws.Range("K33").Value = 296.44446
$ws.Range("L33").Value = 330
$ws.Range("M33").Value = -67.44445999999999
$ws.Range("N33").Value = -788

$ws.Range("H40").Value = 2964
$ws.Range("I40").Value = 2311.6
$ws.Range("J40").Value = 3235.8333
$ws.Range("K40").Value = 2311.6
$ws.Range("L40").Value = 3235.8333
$ws.Range("M40").Value = -2136.6
$ws.Range("N40").Value = -3585.8333

$ws.Range("H51").Value = 3949.6667
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3949.6667
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3949.6667
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4917.6667

$ws.Range("H74").Value = 4253.92
$ws.Range("I74").Value = 3096.2856
$ws.Range("K74").Value = 3096.2856
$ws.Range("M74").Value = -2160.2856

$ws.Range("H77").Value = 4253.92
$ws.Range("I77").Value = 3096.2856
$ws.Range("K77").Value = 15481.428
$ws.Range("M77").Value = -10801.428

$ws.Range("H86").Value = 3066.1614
$ws.Range("I86").Value = 1758.5555
$ws.Range("J86").Value = 4876.6924
$ws.Range("K86").Value = 1758.5555
$ws.Range("L86").Value = 4876.6924
$ws.Range("M86").Value = -635.5554999999999
$ws.Range("N86").Value = -7122.6924

$ws.Range("H89").Value = 3066.1614
$ws.Range("I89").Value = 1758.5555
$ws.Range("J89").Value = 4876.6924
$ws.Range("K89").Value = 8792.7775
$ws.Range("L89").Value = 24383.462
$ws.Range("M89").Value = -3176.7775
$ws.Range("N89").Value = -35615.462

$ws.Range("H116").Value = 3713.625
$ws.Range("I116").Value = 3784.4
$ws.Range("J116").Value = 3595.6667
$ws.Range("K116").Value = 3784.4
$ws.Range("L116").Value = 3595.6667
$ws.Range("M116").Value = -342.4000000000001
$ws.Range("N116").Value = -10479.6667

$ws.Range("H136").Value = 93750
$ws.Range("J136").Value = 93750
$ws.Range("L136").Value = 93750
$ws.Range("N136").Value = -103950

$ws.Range("H138").Value = 6671180
$ws.Range("I138").Value = 1546.4445
$ws.Range("J138").Value = 10422849
$ws.Range("K138").Value = 4639.333500000001
$ws.Range("L138").Value = 31268547
$ws.Range("M138").Value = 500.6664999999994
$ws.Range("N138").Value = -31278827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16141454
$ws.Range("I32").Value = 28577542
$ws.Range("K32").Value = 28577542
$ws.Range("M32").Value = -28577255

$ws.Range("H61").Value = 28575102
$ws.Range("I61").Value = 35716616
$ws.Range("J61").Value = 9057
$ws.Range("K61").Value = 35716616
$ws.Range("L61").Value = 9057
$ws.Range("M61").Value = -35716404
$ws.Range("N61").Value = -9481

$ws.Range("H122").Value = 3150.3447
$ws.Range("I122").Value = 1383.0834
$ws.Range("K122").Value = 4149.2502
$ws.Range("M122").Value = -1699.2502

$ws.Range("I132").Value = 7288.143
$ws.Range("K132").Value = 21864.429
$ws.Range("M132").Value = -19334.429

$ws.Range("H136").Value = 28575102
$ws.Range("I136").Value = 35716616
$ws.Range("J136").Value = 9057
$ws.Range("K136").Value = 107149848
$ws.Range("L136").Value = 27171
$ws.Range("M136").Value = -107147298
$ws.Range("N136").Value = -32271

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H94").Value = 2229.3225
$ws.Range("I94").Value = 1942.0416
$ws.Range("J94").Value = 3214.2856
$ws.Range("K94").Value = 1942.0416
$ws.Range("L94").Value = 3214.2856
$ws.Range("M94").Value = -1491.0416
$ws.Range("N94").Value = -4116.2856

$ws.Range("H105").Value = 6220.15
$ws.Range("I105").Value = 7219.8125
$ws.Range("K105").Value = 7219.8125
$ws.Range("M105").Value = -5472.8125

$ws.Range("H140").Value = 176489.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 176489.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 176489.8
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -186849.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 46892.176
$ws.Range("I132").Value = 66295.06
$ws.Range("K132").Value = 198885.18
$ws.Range("M132").Value = -196355.18

$ws.Range("H134").Value = 1452.1666
$ws.Range("I134").Value = 1022.3684
$ws.Range("K134").Value = 3067.1052
$ws.Range("M134").Value = -532.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H33").Value = 95.5
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 126.666664
$ws.Range("K33").Value = 12
$ws.Range("L33").Value = 759.999984
$ws.Range("M33").Value = 271
$ws.Range("N33").Value = -1325.999984

$ws.Range("H80").Value = 2625.75
$ws.Range("J80").Value = 2625.75
$ws.Range("L80").Value = 7877.25
$ws.Range("N80").Value = -9749.25

$ws.Range("H83").Value = 2625.75
$ws.Range("J83").Value = 2625.75
$ws.Range("L83").Value = 23631.75
$ws.Range("N83").Value = -32991.75

$ws.Range("H122").Value = 1539.2222
$ws.Range("I122").Value = 386.75
$ws.Range("J122").Value = 2461.2
$ws.Range("K122").Value = 3480.75
$ws.Range("L122").Value = 22150.8
$ws.Range("M122").Value = -1030.75
$ws.Range("N122").Value = -27050.8

$ws.Range("H123").Value = 2800
$ws.Range("I123").Value = 2800
$ws.Range("K123").Value = 8400
$ws.Range("M123").Value = -5950

$ws.Range("H131").Value = 30325.615
$ws.Range("I131").Value = 58658.61
$ws.Range("J131").Value = 6040.1904
$ws.Range("K131").Value = 175975.83
$ws.Range("L131").Value = 18120.5712
$ws.Range("M131").Value = -170935.83
$ws.Range("N131").Value = -28200.5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3094.7334
$ws.Range("J80").Value = 3555.6667
$ws.Range("L80").Value = 3555.6667
$ws.Range("N80").Value = -5551.6667

$ws.Range("H83").Value = 3094.7334
$ws.Range("J83").Value = 3555.6667
$ws.Range("L83").Value = 17778.3335
$ws.Range("N83").Value = -27762.3335

$ws.Range("H97").Value = 2140.25
$ws.Range("I97").Value = 676.375
$ws.Range("K97").Value = 676.375
$ws.Range("M97").Value = -180.375

$ws.Range("H122").Value = 2004.2
$ws.Range("I122").Value = 1825.7368
$ws.Range("K122").Value = 5477.2104
$ws.Range("M122").Value = -3027.2104

$ws.Range("H126").Value = 25008750
$ws.Range("I126").Value = 16678332
$ws.Range("J126").Value = 33339166
$ws.Range("K126").Value = 50034996
$ws.Range("L126").Value = 100017498
$ws.Range("M126").Value = -50032526
$ws.Range("N126").Value = -100022438

$ws.Range("H132").Value = 2555.04
$ws.Range("I132").Value = 2465.6
$ws.Range("J132").Value = 2912.8
$ws.Range("K132").Value = 7396.799999999999
$ws.Range("L132").Value = 8738.400000000001
$ws.Range("M132").Value = -4866.799999999999
$ws.Range("N132").Value = -13798.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1604
$ws.Range("I46").Value = 952.6486
$ws.Range("J46").Value = 3325.4285
$ws.Range("K46").Value = 952.6486
$ws.Range("L46").Value = 3325.4285
$ws.Range("M46").Value = -764.6486
$ws.Range("N46").Value = -3701.4285

$ws.Range("H56").Value = 22749.75
$ws.Range("I56").Value = 22749.75
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 22749.75
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -22058.75
$ws.Range("N56").ClearContents()

$ws.Range("H93").Value = 1349.9546
$ws.Range("I93").Value = 968.75
$ws.Range("J93").Value = 2366.5
$ws.Range("K93").Value = 968.75
$ws.Range("L93").Value = 2366.5
$ws.Range("M93").Value = 279.25
$ws.Range("N93").Value = -4862.5

$ws.Range("H132").Value = 125002060
$ws.Range("I132").Value = 1915.1428
$ws.Range("J132").Value = 222224400
$ws.Range("K132").Value = 5745.428400000001
$ws.Range("L132").Value = 666673200
$ws.Range("M132").Value = -3215.428400000001
$ws.Range("N132").Value = -666678260

$ws.Range("H139").Value = 55000
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6984.75
$ws.Range("J74").Value = 7268.2856
$ws.Range("L74").Value = 7268.2856
$ws.Range("N74").Value = -9140.285599999999

$ws.Range("H77").Value = 6984.75
$ws.Range("J77").Value = 7268.2856
$ws.Range("L77").Value = 21804.8568
$ws.Range("N77").Value = -31164.8568

$ws.Range("H86").Value = 20041458
$ws.Range("J86").Value = 59999
$ws.Range("L86").Value = 59999
$ws.Range("N86").Value = -62245

$ws.Range("H89").Value = 20041458
$ws.Range("J89").Value = 59999
$ws.Range("L89").Value = 299995
$ws.Range("N89").Value = -311227

$ws.Range("H132").Value = 6386.057
$ws.Range("I132").Value = 6349
$ws.Range("K132").Value = 19047
$ws.Range("M132").Value = -16517
